$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 43
$ws.Range("H43").Value = 785.05554
$ws.Range("I43").Value = 487.66666
$ws.Range("J43").Value = 933.75
$ws.Range("K43").Value = 487.66666
$ws.Range("L43").Value = 933.75
$ws.Range("M43").Value = -418.66666
$ws.Range("N43").Value = -1071.75
# Row 86
$ws.Range("H86").Value = 73762.07000000001
$ws.Range("I86").Value = 85349.664
$ws.Range("J86").Value = 4236.5
$ws.Range("K86").Value = 85349.664
$ws.Range("L86").Value = 4236.5
$ws.Range("M86").Value = -84226.664
$ws.Range("N86").Value = -6482.5
# Row 89
$ws.Range("H89").Value = 73762.07000000001
$ws.Range("I89").Value = 85349.664
$ws.Range("J89").Value = 4236.5
$ws.Range("K89").Value = 426748.32
$ws.Range("L89").Value = 21182.5
$ws.Range("M89").Value = -421132.32
$ws.Range("N89").Value = -32414.5
# Row 113
$ws.Range("H113").Value = 3592
$ws.Range("I113").Value = 2200
$ws.Range("J113").Value = 3790.8572
$ws.Range("K113").Value = 2200
$ws.Range("L113").Value = 3790.8572
$ws.Range("M113").Value = 1054
$ws.Range("N113").Value = -10298.8572

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1493.3214
$ws.Range("I45").Value = 1517.25
$ws.Range("J45").Value = 1349.75
$ws.Range("K45").Value = 1517.25
$ws.Range("L45").Value = 1349.75
$ws.Range("M45").Value = -1140.25
$ws.Range("N45").Value = -2103.75
# Row 102
$ws.Range("H102").Value = 3670
$ws.Range("I102").Value = 3601.818
$ws.Range("J102").Value = 3820
$ws.Range("K102").Value = 3601.818
$ws.Range("L102").Value = 3820
$ws.Range("M102").Value = -1979.818
$ws.Range("N102").Value = -7064
# Row 122
$ws.Range("H122").Value = 1921.0952
$ws.Range("I122").Value = 1768
$ws.Range("J122").Value = 2125.2222
$ws.Range("K122").Value = 5304
$ws.Range("L122").Value = 6375.6666
$ws.Range("M122").Value = -2854
$ws.Range("N122").Value = -11275.6666
# Row 132
$ws.Range("H132").Value = 2876.9285
$ws.Range("I132").Value = 2267.1765
$ws.Range("J132").Value = 3819.2727
$ws.Range("K132").Value = 6801.529500000001
$ws.Range("L132").Value = 11457.8181
$ws.Range("M132").Value = -4271.529500000001
$ws.Range("N132").Value = -16517.8181

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 1611.1694
$ws.Range("I86").Value = 1629.8462
$ws.Range("J86").Value = 1472.4286
$ws.Range("K86").Value = 1629.8462
$ws.Range("L86").Value = 1472.4286
$ws.Range("M86").Value = -506.8462
$ws.Range("N86").Value = -3718.4286
# Row 89
$ws.Range("H89").Value = 1611.1694
$ws.Range("I89").Value = 1629.8462
$ws.Range("J89").Value = 1472.4286
$ws.Range("K89").Value = 8149.231
$ws.Range("L89").Value = 7362.143
$ws.Range("M89").Value = -2533.231
$ws.Range("N89").Value = -18594.143
# Row 103
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
# Row 107
$ws.Range("H107").Value = 1256.7778
$ws.Range("I107").Value = 1102.2
$ws.Range("J107").Value = 1450
$ws.Range("K107").Value = 1102.2
$ws.Range("L107").Value = 1450
$ws.Range("M107").Value = 817.8
$ws.Range("N107").Value = -5290

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 122
$ws.Range("H122").Value = 9960.333000000001
$ws.Range("I122").Value = 5462.4614
$ws.Range("J122").Value = 21654.8
$ws.Range("K122").Value = 16387.3842
$ws.Range("L122").Value = 64964.39999999999
$ws.Range("M122").Value = -13937.3842
$ws.Range("N122").Value = -69864.39999999999
# Row 132
$ws.Range("H132").Value = 2706.2666
$ws.Range("I132").Value = 2602.64
$ws.Range("J132").Value = 2835.8
$ws.Range("K132").Value = 7807.92
$ws.Range("L132").Value = 8507.400000000001
$ws.Range("M132").Value = -5277.92
$ws.Range("N132").Value = -13567.4

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 20
$ws.Range("H20").Value = 1810.1
$ws.Range("J20").Value = 2800
$ws.Range("L20").Value = 8400
$ws.Range("N20").Value = -8854

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 13749
$ws.Range("I122").Value = 25498
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 76494
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -74044
$ws.Range("N122").Value = -10900
# Row 126
$ws.Range("H126").Value = 2361.9355
$ws.Range("I126").Value = 1689.4117
$ws.Range("J126").Value = 3178.5715
$ws.Range("K126").Value = 5068.2351
$ws.Range("L126").Value = 9535.7145
$ws.Range("M126").Value = -2598.2351
$ws.Range("N126").Value = -14475.7145
# Row 132
$ws.Range("H132").Value = 6321.607
$ws.Range("I132").Value = 2410.1365
$ws.Range("J132").Value = 20663.666
$ws.Range("K132").Value = 7230.4095
$ws.Range("L132").Value = 61990.99800000001
$ws.Range("M132").Value = -4700.4095
$ws.Range("N132").Value = -67050.99800000001
# Row 137
$ws.Range("H137").Value = 49800
$ws.Range("J137").Value = 49800
$ws.Range("L137").Value = 49800
$ws.Range("N137").Value = -60000

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 812.5333000000001
$ws.Range("I22").Value = 789.2
$ws.Range("J22").Value = 824.2
$ws.Range("K22").Value = 789.2
$ws.Range("L22").Value = 824.2
$ws.Range("M22").Value = -494.2
$ws.Range("N22").Value = -1414.2
# Row 27
$ws.Range("H27").Value = 812.5333000000001
$ws.Range("I27").Value = 789.2
$ws.Range("J27").Value = 824.2
$ws.Range("K27").Value = 789.2
$ws.Range("L27").Value = 824.2
$ws.Range("M27").Value = -682.2
$ws.Range("N27").Value = -1038.2
# Row 68
$ws.Range("H68").Value = 3088.889
$ws.Range("I68").Value = 2614.2856
$ws.Range("J68").Value = 4750
$ws.Range("K68").Value = 2614.2856
$ws.Range("L68").Value = 4750
$ws.Range("M68").Value = -1865.2856
$ws.Range("N68").Value = -6248
# Row 71
$ws.Range("H71").Value = 3088.889
$ws.Range("I71").Value = 2614.2856
$ws.Range("J71").Value = 4750
$ws.Range("K71").Value = 13071.428
$ws.Range("L71").Value = 23750
$ws.Range("M71").Value = -9327.428
$ws.Range("N71").Value = -31238
# Row 100
$ws.Range("H100").Value = 4171.7856
$ws.Range("I100").Value = 3700.4167
$ws.Range("K100").Value = 3700.4167
$ws.Range("M100").Value = -3159.4167
# Row 122
$ws.Range("H122").Value = 6718.383
$ws.Range("I122").Value = 6641.3076
$ws.Range("K122").Value = 19923.9228
$ws.Range("M122").Value = -17473.9228
# Row 132
$ws.Range("H132").Value = 4755.7827
$ws.Range("I132").Value = 5057.0835
$ws.Range("J132").Value = 4427.091
$ws.Range("K132").Value = 15171.2505
$ws.Range("L132").Value = 13281.273
$ws.Range("M132").Value = -12641.2505
$ws.Range("N132").Value = -18341.273
# Row 136
$ws.Range("H136").Value = 5676.75
$ws.Range("I136").Value = 3407.15
$ws.Range("J136").Value = 8513.75
$ws.Range("K136").Value = 10221.45
$ws.Range("L136").Value = 25541.25
$ws.Range("M136").Value = -7671.450000000001
$ws.Range("N136").Value = -30641.25

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 28574300
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 28574300
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 57148600
$ws.Range("N81").Value = -57150722
$ws.Range("M81").ClearContents()
# Row 84
$ws.Range("H84").Value = 28574300
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 28574300
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 285743000
$ws.Range("N84").Value = -285753608
$ws.Range("M84").ClearContents()
# Row 107
$ws.Range("H107").Value = 2294.4211
$ws.Range("I107").Value = 516.6667
$ws.Range("J107").Value = 3894.4
$ws.Range("K107").Value = 1550.0001
$ws.Range("L107").Value = 11683.2
$ws.Range("M107").Value = 369.9999
$ws.Range("N107").Value = -15523.2
# Row 122
$ws.Range("H122").Value = 16601.428
$ws.Range("I122").Value = 1200
$ws.Range("K122").Value = 3600
$ws.Range("M122").Value = -1150
# Row 126
$ws.Range("H126").Value = 2100.4
$ws.Range("I126").Value = 2063
$ws.Range("J126").Value = 2250
$ws.Range("K126").Value = 6189
$ws.Range("L126").Value = 6750
$ws.Range("M126").Value = -3719
$ws.Range("N126").Value = -11690
# Row 132
$ws.Range("H132").Value = 1819.186
$ws.Range("I132").Value = 1110.5862
$ws.Range("J132").Value = 3287
$ws.Range("K132").Value = 3331.7586
$ws.Range("L132").Value = 9861
$ws.Range("M132").Value = -801.7586000000001
$ws.Range("N132").Value = -14921
# Row 133
$ws.Range("H133").Value = 40905
$ws.Range("J133").Value = 40905
$ws.Range("L133").Value = 40905
$ws.Range("N133").Value = -51025
# Row 136
$ws.Range("H136").Value = 5905.3784
$ws.Range("I136").Value = 1996.45
$ws.Range("K136").Value = 5989.35
$ws.Range("M136").Value = -3439.35
